$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.328.93"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.089.56"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.93%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "386.51"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.40"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("E7").Value = "  -1.70%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  -1.59%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.85"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.26%  "

$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.579.01"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.47"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.75"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.102.12"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.988"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.67"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.431.88"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.21"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.30%  "

$ws.Range("E21").Value = "  -0.34%  "

$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.95"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "265.45"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.14"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.99"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.29"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.03%  "

$ws.Range("E28").Value = "  -6.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("E30").Value = "  -5.05%  "

$ws.Range("E31").Value = "  -2.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.39"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.57"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.35%  "

$ws.Range("E34").Value = "  +5.64%  "

$ws.Range("E35").Value = "  -1.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "49.92"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.18%  "

$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.35"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.90%  "

$ws.Range("E39").Value = "  -1.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "129.43"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.81%  "

$ws.Range("E41").Value = "  -0.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.115"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.48"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.47%  "

$ws.Range("E44").Value = "  +0.74%  "

$ws.Range("E45").Value = "  -3.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.89"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.50"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.75%  "

$ws.Range("E48").Value = "  -0.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.069.54"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.926"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +17.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0327"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.77%  "
